# Update the header date line (unique text in the document).
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-01-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-22 Monday", 2)

# Update the division problems in the single table.
# Cell text is set directly (rather than via Find/Replace) because some
# problems share identical text (e.g. two cells both contain "95÷6="
# but must become different values), and Find/Replace operates on the
# whole document rather than being confined to a single cell's range.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "43÷8="
$t.Cell(1,2).Range.Text = "64÷8="
$t.Cell(1,3).Range.Text = "89÷8="
$t.Cell(1,4).Range.Text = "72÷5="
$t.Cell(1,5).Range.Text = "54÷8="

$t.Cell(5,1).Range.Text = "26÷3="
$t.Cell(5,2).Range.Text = "65÷2="
$t.Cell(5,3).Range.Text = "93÷6="
$t.Cell(5,4).Range.Text = "89÷4="
$t.Cell(5,5).Range.Text = "25÷9="

$t.Cell(9,1).Range.Text = "75÷2="
$t.Cell(9,2).Range.Text = "37÷5="
$t.Cell(9,3).Range.Text = "37÷8="
$t.Cell(9,4).Range.Text = "67÷2="
$t.Cell(9,5).Range.Text = "15÷5="

$t.Cell(13,1).Range.Text = "85÷4="
$t.Cell(13,2).Range.Text = "24÷6="
$t.Cell(13,3).Range.Text = "24÷9="
$t.Cell(13,4).Range.Text = "26÷6="
$t.Cell(13,5).Range.Text = "28÷4="

$t.Cell(17,1).Range.Text = "89÷4="
$t.Cell(17,2).Range.Text = "85÷8="
$t.Cell(17,3).Range.Text = "75÷5="
$t.Cell(17,4).Range.Text = "63÷7="
$t.Cell(17,5).Range.Text = "80÷2="
